$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row {
    param(
        [int]$Row1,
        [int]$Row2
    )
    # Swap values in columns B (2) through AD (30) between two rows
    for ($col = 2; $col -le 30; $col++) {
        $cell1 = $ws.Cells.Item($Row1, $col)
        $cell2 = $ws.Cells.Item($Row2, $col)
        $v1 = $cell1.Value2
        $v2 = $cell2.Value2
        $cell1.Value2 = $v2
        $cell2.Value2 = $v1
    }
}

Swap-Row 45 46
Swap-Row 54 55
Swap-Row 78 79
Swap-Row 81 82
Swap-Row 83 84
Swap-Row 87 88
Swap-Row 108 109
Swap-Row 118 119
Swap-Row 122 123
Swap-Row 124 125
Swap-Row 129 130
